# Sam's "working model sans random vars" edit:
# Clears out the random-variable rows (Food name + the 3 numeric columns)
# so the sheet holds a "blank template" for those rows, widens column A to
# fit the remaining food names, and leaves the selection/scroll position
# where Sam had left off working (around row 20-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Food/Refrigerated/StoreLoss/HomeWaste values get cleared
# (style/formatting is retained - only the values are removed).
$rowsToClear = @(14, 21, 22, 27, 30, 45, 47)
foreach ($r in $rowsToClear) {
    $ws.Range("A$r`:D$r").ClearContents()
}

# Widen column A to fit the (now shorter) set of longest food names.
$ws.Columns("A:A").ColumnWidth = 20.33

# Restore the working selection/scroll state.
$win = $excel.ActiveWindow
$win.ScrollRow = 20
$win.ScrollColumn = 1
$ws.Range("E28").Select() | Out-Null
